# Update placeholder survey dates to zero-padded month/day format.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

$ws.Range("B38").Value = "????-03-11"
$ws.Range("B39").Value = "????-04-08"
$ws.Range("B40").Value = "????-05-06"
$ws.Range("B41").Value = "????-05-18"
$ws.Range("B43").Value = "????-07-08"
$ws.Range("B44").Value = "????-07-08"
$ws.Range("B42").Value = "????-06-09"

# Move the active selection to match the saved view state.
$ws.Range("I34").Select() | Out-Null
